$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new blank row above the old row 13 ("Programa resumido:" row),
# shifting all subsequent rows down by one. Row heights carry over correctly
# (the newly inserted row gets no custom height, matching the target).
$ws.Rows.Item(13).Insert()

# The inserted row copied formatting from the row above into A13; the target
# layout has no value in column A for this row, so drop that stray cell
# entirely (Clear removes both content and formatting so no empty <c> tag
# is serialized).
$ws.Range("A13").Clear()

# Populate the now-empty B13/C13 with the correct text and adopt the same
# style used elsewhere in columns B/C (copy format from a known-good cell).
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# --- Fix up the text content for every row from 10 through 24 ---

$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "Introduzir as técnicas de fabricação de dispositivos e circuitos integrados em microeletrônica. Apresenta os princípios, técnicas, equipamentos e softwares utilizados na simulação e fabricação de dispositivos em silício e arseneto de gálio de uma maneira global e genérica."
$ws.Range("C10").Value = "Introduzir as técnicas de fabricação de dispositivos e circuitos integrados em microeletrônica. Apresenta os princípios, técnicas, equipamentos e softwares utilizados na simulação e fabricação de dispositivos em silício e arseneto de gálio de uma maneira global e genérica."

$ws.Range("A11").Value = "Objectives:"

$ws.Range("A12").Value = "Docentes responsáveis:"

$ws.Range("B13").Value = "5840535 - Messias Borges Silva"
$ws.Range("C13").Value = "5840535 - Messias Borges Silva"

$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "Processos de fabricação típicos e de montagem de sistemas eletroeletronicos"
$ws.Range("C14").Value = "Processos de fabricação típicos e de montagem de sistemas eletroeletronicos"

$ws.Range("A15").Value = "Short syllabus:"

$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "Técnicas de obtenção de silício cristalino e arseneto de gálio. Processamento de lâminas de silício e GaAs. Rede e Defeitos cristalinos. Oxidação. Dopagem (difusão e implantação iônica). Litografia: (Fabricação de Fotomáscaras; Gerador de Padrões). Epitaxia. Deposição em Fase Vapor (CVD. PECVD e LPCVD). Decapagem (úmida e seca). Equipamentos e técnicas de metalização (Evaporação térmica, feixe iônico, bombardeamento catódico (`"sputtering`"), Caracterização de etapas de processo de fabricação. Simulação de processos de fabricação. Montagem de Sistemas Eletro Eletrônicos"
$ws.Range("C16").Value = "Técnicas de obtenção de silício cristalino e arseneto de gálio. Processamento de lâminas de silício e GaAs. Rede e Defeitos cristalinos. Oxidação. Dopagem (difusão e implantação iônica). Litografia: (Fabricação de Fotomáscaras; Gerador de Padrões). Epitaxia. Deposição em Fase Vapor (CVD. PECVD e LPCVD). Decapagem (úmida e seca). Equipamentos e técnicas de metalização (Evaporação térmica, feixe iônico, bombardeamento catódico (`"sputtering`"), Caracterização de etapas de processo de fabricação. Simulação de processos de fabricação. Montagem de Sistemas Eletro Eletrônicos"

$ws.Range("A17").Value = "Syllabus:"

$ws.Range("A18").Value = "Avaliação:"

$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "Aulas expositivas auxiliadas por transparências e métodos multimídia (datashow + computador) São previstas três aulas práticas com visitas em ambientes de fabricação de circuitos integrados (salas limpas) e montagem de produtos eletroeletrônicos. Duas prvas e relatórios de trabalhos práticos e visitas."
$ws.Range("C19").Value = "Aulas expositivas auxiliadas por transparências e métodos multimídia (datashow + computador) São previstas três aulas práticas com visitas em ambientes de fabricação de circuitos integrados (salas limpas) e montagem de produtos eletroeletrônicos. Duas prvas e relatórios de trabalhos práticos e visitas."

$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "Média geral MG = [0,6 x (média aritmética das 2 provas feitas) + 0,2 x média dos trabalhos práticos + 0,2 x (média aritmética dos relatórios de visitas)] >= 5,0. Os testes serão dados nos 10 minutos finais de N aulas escolhidas aleatoriamente. A prova substitutiva substitui a prova em que o aluno faltou."
$ws.Range("C20").Value = "Média geral MG = [0,6 x (média aritmética das 2 provas feitas) + 0,2 x média dos trabalhos práticos + 0,2 x (média aritmética dos relatórios de visitas)] >= 5,0. Os testes serão dados nos 10 minutos finais de N aulas escolhidas aleatoriamente. A prova substitutiva substitui a prova em que o aluno faltou."

$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "1 (uma) prova de recuperação."
$ws.Range("C21").Value = "1 (uma) prova de recuperação."

$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "01 Stephen A. Campbell, The Science and Engineering of Microelectronic Fabrication, 2nd edition, Oxford University Press, 2001;02 S.M. Sze - VLSI Technology, McGraw-Hill, 1988;03 V. Baranauskas, org., Processos em Microeletrônica, SBV e SBMicro, 1990;04 R. A. Levy, Microelectronic Materials and Processes, Kluwer Academic Publ., 1989;[05] M. Madou, Fundamentals of Microfabrication, CRC press, 1997."
$ws.Range("C22").Value = "01 Stephen A. Campbell, The Science and Engineering of Microelectronic Fabrication, 2nd edition, Oxford University Press, 2001;02 S.M. Sze - VLSI Technology, McGraw-Hill, 1988;03 V. Baranauskas, org., Processos em Microeletrônica, SBV e SBMicro, 1990;04 R. A. Levy, Microelectronic Materials and Processes, Kluwer Academic Publ., 1989;[05] M. Madou, Fundamentals of Microfabrication, CRC press, 1997."

$ws.Range("A23").Value = "Requisitos:"

$ws.Range("B24").Value = "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)`n"

# --- Column layout cleanup: column A's width/style entry used to overlap
# with column B's in <cols> (min="1" max="2"). Touching column B's width
# forces the engine to split that merged range into per-column entries so
# column A ends up on its own, matching the target layout.
$ws.Columns.Item(2).ColumnWidth = 60.7109375
